# Quarterly dictionary workbook — "Categories" sheet:
# the isMissing column (C) was stored as the literal text "FALSE" for every
# row; convert it to an actual boolean FALSE value (commit: "change text
# value to boolean value quarterly variables").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Rows 2-29 hold data; column C ("isMissing") is text "FALSE" everywhere.
$ws.Range("C2:C29").Value = $false

# Reflect the author's UI focus at save time: Categories tab active with
# C3:C29 selected.
$ws.Activate()
$ws.Range("C3:C29").Select() | Out-Null
